$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns
$ws.Range("D1").Value = "CompletedExercises"
$ws.Range("E1").Value = "TestResults"

# Existing user "Dima" (row 2) becomes "Dima1" and gets progress data
$ws.Range("B2").Value = "Dima1"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "4/20"

# New user "Dima" added as row 3, with the same password as row 2
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Dima"
$ws.Range("C2").Copy($ws.Range("C3"))
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "3/20"
